# Observations-summary: add the "US Core Average Blood Pressure Profile"
# rows (USCDI4 gh-pages deploy) as rows 2-4 of the Observations sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the three new rows the same cell formatting (border/wrap/valign) that
# the header row already uses (style index 2 in the original workbook) by
# copying the header row's format down into A2:K4.
$ws.Range("A1:K1").Copy() | Out-Null
$ws.Range("A2:K4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2: us-core-average-blood-pressure profile header row
$ws.Cells.Item(2, 1).Value = 'us-core-average-blood-pressure'
$ws.Cells.Item(2, 2).Value = 'US Core Average Blood Pressure Profile'
$ws.Cells.Item(2, 3).Value = 'null#vital-signs'
$ws.Cells.Item(2, 5).Value = 'LOINC#96607-7'
$ws.Cells.Item(2, 7).Value = 'dateTimeĵ, Periodĵ'
$ws.Cells.Item(2, 8).Value = 'Quantityĵ, CodeableConceptĵ, stringĵ, booleanĵ, integerĵ, Rangeĵ, Ratioĵ, SampledDataĵ, timeĵ, dateTimeĵ, Periodĵ'
$ws.Cells.Item(2, 9).Value = 'optional'

# Row 3: second component code for the same profile
$ws.Cells.Item(3, 2).Value = 'US Core Average Blood Pressure Profile'
$ws.Cells.Item(3, 5).Value = 'LOINC#96608-5'
$ws.Cells.Item(3, 8).Value = 'Quantityĵ'
$ws.Cells.Item(3, 9).Value = 'optional'

# Row 4: third component code for the same profile
$ws.Cells.Item(4, 2).Value = 'US Core Average Blood Pressure Profile'
$ws.Cells.Item(4, 5).Value = 'LOINC#96609-3'
$ws.Cells.Item(4, 8).Value = 'Quantityĵ'
$ws.Cells.Item(4, 9).Value = 'optional'
